# Adds two new worksheets ("existingAccount" and "notExistingAccount") with
# login/password test data, mirroring the commit "Added data to Excel".

$wb = $excel.ActiveWorkbook

# --- Create the two new worksheets, after the existing ones -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsExisting = $wb.Worksheets.Add($null, $lastSheet)
$wsExisting.Name = "existingAccount"

$wsNotExisting = $wb.Worksheets.Add($null, $wsExisting)
$wsNotExisting.Name = "notExistingAccount"

# --- existingAccount sheet ---------------------------------------------------
$wsExisting.Range("A1").Value = "jan@o2.pl"
$wsExisting.Range("A1").Font.Color = 3880497   # RGB(0x31,0x36,0x3B) -> FF31363B
$wsExisting.Range("B1").Value = "haslo"

# --- notExistingAccount sheet ------------------------------------------------
$wsNotExisting.Range("A1").Value = "jann@o2.pl"
$wsNotExisting.Range("B1").Value = "haslo1"

$wsNotExisting.Range("A2").Value = "jan"
$wsNotExisting.Range("B2").Value = "haslo1"

$wsNotExisting.Range("A3").Value = "jan@o2.pl"
$wsNotExisting.Range("B3").Value = "haslo1"

$wsNotExisting.Range("A4").Value = "jann@o2.pl"
$wsNotExisting.Range("B4").Value = "haslo"

# Hyperlinks (mailto:) on the e-mail cells, re-using the blue Arial style that
# is already used for the hyperlinks on the other two sheets.
$wsNotExisting.Hyperlinks.Add($wsNotExisting.Range("A1"), "mailto:jann@o2.pl", "", "", "jann@o2.pl")
$wsNotExisting.Range("A1").Font.Name = "Arial"
$wsNotExisting.Range("A1").Font.Underline = -4142
$wsNotExisting.Range("A1").Font.Color = 16711680

$wsNotExisting.Hyperlinks.Add($wsNotExisting.Range("A3"), "mailto:jan@o2.pl", "", "", "jan@o2.pl")
$wsNotExisting.Range("A3").Font.Name = "Arial"
$wsNotExisting.Range("A3").Font.Underline = -4142
$wsNotExisting.Range("A3").Font.Color = 16711680

$wsNotExisting.Hyperlinks.Add($wsNotExisting.Range("A4"), "mailto:jann@o2.pl", "", "", "jann@o2.pl")
$wsNotExisting.Range("A4").Font.Name = "Arial"
$wsNotExisting.Range("A4").Font.Underline = -4142
$wsNotExisting.Range("A4").Font.Color = 16711680

# Restore the cursor positions left behind on the two new sheets.
$wsExisting.Range("B12").Select() | Out-Null
$wsNotExisting.Range("B6").Select() | Out-Null

# Make "notExistingAccount" the active tab (as in the target workbook).
$wsNotExisting.Activate()
